$d = $word.ActiveDocument

# Locate the anchor paragraph: "Order. Comparison. Relations. Upper Ontology
# assertions. Augmentations. TBD." -- the new block of paragraphs is inserted
# immediately after it (and before the pre-existing blank paragraph that leads
# into "Functional API: Monads / Transforms").
$idx = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Order. Comparison. Relations. Upper Ontology assertions. Augmentations. TBD.*") {
        $idx = $i
        break
    }
}

if ($idx -eq 0) {
    throw "Anchor paragraph not found"
}

# The new content: blank separator lines plus the labelled lines describing the
# Relation / Relationship (Tabular / OGM) section. Only the heading line is
# underlined.
$items = @(
    @{ Text = $null; Underline = $false },
    @{ Text = "Relation / Relationship: Tabular / OGM (Object Graph Mapper):"; Underline = $true },
    @{ Text = $null; Underline = $false },
    @{ Text = "I/O: (Class, ClassID, Attribute, Value);"; Underline = $false },
    @{ Text = $null; Underline = $false },
    @{ Text = "Class: Table / Object Type."; Underline = $false },
    @{ Text = "ClassID: PK / Object ID."; Underline = $false },
    @{ Text = "Attribute: Column / Member."; Underline = $false },
    @{ Text = "Value: Cell / Field Value."; Underline = $false },
    @{ Text = $null; Underline = $false },
    @{ Text = "Subject Kind: Relation / Domain."; Underline = $false },
    @{ Text = $null; Underline = $false },
    @{ Text = "Predicate Kind: Relationship."; Underline = $false },
    @{ Text = $null; Underline = $false },
    @{ Text = "Object Kind: Mapping / Range."; Underline = $false },
    @{ Text = $null; Underline = $false },
    @{ Text = "Dataflow: Reactive Functional Augmentation / Integration APIs."; Underline = $false },
    @{ Text = $null; Underline = $false },
    @{ Text = "Indices: Apply functional mappings expansion."; Underline = $false }
)

# First pass: create all the new (plain) paragraphs and fill in their text.
# Underline is applied afterwards in a second pass so that it never "bleeds"
# forward into paragraphs inserted later (InsertParagraphAfter clones the
# formatting of the paragraph it splits from).
$underlineIdx = @()
foreach ($item in $items) {
    $r = $d.Paragraphs.Item($idx).Range
    $r.InsertParagraphAfter()
    $idx = $idx + 1

    if ($item.Text -ne $null) {
        $nr = $d.Paragraphs.Item($idx).Range
        $nr.InsertBefore($item.Text)
    }

    if ($item.Underline) {
        $underlineIdx += $idx
    }
}

# Second pass: apply the underline formatting to just the heading paragraph.
foreach ($ui in $underlineIdx) {
    $d.Paragraphs.Item($ui).Range.Font.Underline = 1
}
